$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 9,20
$data[0,0] = "ECs"
$data[0,1] = "Col1a2"
$data[0,2] = "Itga11"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 10.45491533333333
$data[0,7] = 31.364746
$data[0,8] = 0.0134573334963438
$data[0,9] = 0.0134573334963438
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.04023833333333333
$data[0,13] = 0.120715
$data[0,14] = 0.001153919673903629
$data[0,15] = 0.001153919673903629
$data[0,16] = 0.4206883681544445
$data[0,17] = 3.78619531339
$data[0,18] = 0.00001552868187971342
$data[0,19] = 0.00001552868187971342
$data[1,0] = "ECs"
$data[1,1] = "Col1a2"
$data[1,2] = "Itga11"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 10.45491533333333
$data[1,7] = 31.364746
$data[1,8] = 0.0134573334963438
$data[1,9] = 0.0134573334963438
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 33.94639966666666
$data[1,13] = 101.839199
$data[1,14] = 0.9734851120464462
$data[1,15] = 0.9734851120464463
$data[1,16] = 354.9067343864949
$data[1,17] = 3194.160609478454
$data[1,18] = 0.01310051380653464
$data[1,19] = 0.01310051380653464
$data[2,0] = "ECs"
$data[2,1] = "Col1a2"
$data[2,2] = "Itga11"
$data[2,3] = "sCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 10.45491533333333
$data[2,7] = 31.364746
$data[2,8] = 0.0134573334963438
$data[2,9] = 0.0134573334963438
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.8843623333333334
$data[2,13] = 2.653087
$data[2,14] = 0.02536096827965006
$data[2,15] = 0.02536096827965006
$data[2,16] = 9.245933318989112
$data[2,17] = 83.21339987090201
$data[2,18] = 0.0003412910079294474
$data[2,19] = 0.0003412910079294474
$data[3,0] = "FAPs"
$data[3,1] = "Col1a2"
$data[3,2] = "Itga11"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 735.4993083333334
$data[3,7] = 2206.497925
$data[3,8] = 0.9467182815928301
$data[3,9] = 0.9467182815928301
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.04023833333333333
$data[3,13] = 0.120715
$data[3,14] = 0.001153919673903629
$data[3,15] = 0.001153919673903629
$data[3,16] = 29.59526633515278
$data[3,17] = 266.357397016375
$data[3,18] = 0.001092436850774202
$data[3,19] = 0.001092436850774202
$data[4,0] = "FAPs"
$data[4,1] = "Col1a2"
$data[4,2] = "Itga11"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 735.4993083333334
$data[4,7] = 2206.497925
$data[4,8] = 0.9467182815928301
$data[4,9] = 0.9467182815928301
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 33.94639966666666
$data[4,13] = 101.839199
$data[4,14] = 0.9734851120464462
$data[4,15] = 0.9734851120464463
$data[4,16] = 24967.55347524023
$data[4,17] = 224707.9812771621
$data[4,18] = 0.9216161524328152
$data[4,19] = 0.9216161524328152
$data[5,0] = "FAPs"
$data[5,1] = "Col1a2"
$data[5,2] = "Itga11"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 735.4993083333334
$data[5,7] = 2206.497925
$data[5,8] = 0.9467182815928301
$data[5,9] = 0.9467182815928301
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.8843623333333334
$data[5,13] = 2.653087
$data[5,14] = 0.02536096827965006
$data[5,15] = 0.02536096827965006
$data[5,16] = 650.4478844827195
$data[5,17] = 5854.030960344476
$data[5,18] = 0.02400969230924058
$data[5,19] = 0.02400969230924058
$data[6,0] = "sCs"
$data[6,1] = "Col1a2"
$data[6,2] = "Itga11"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 30.939307
$data[6,7] = 92.81792100000001
$data[6,8] = 0.03982438491082609
$data[6,9] = 0.03982438491082609
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.04023833333333333
$data[6,13] = 0.120715
$data[6,14] = 0.001153919673903629
$data[6,15] = 0.001153919673903629
$data[6,16] = 1.244946148168333
$data[6,17] = 11.204515333515
$data[6,18] = 0.00004595414124971304
$data[6,19] = 0.00004595414124971304
$data[7,0] = "sCs"
$data[7,1] = "Col1a2"
$data[7,2] = "Itga11"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 30.939307
$data[7,7] = 92.81792100000001
$data[7,8] = 0.03982438491082609
$data[7,9] = 0.03982438491082609
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 33.94639966666666
$data[7,13] = 101.839199
$data[7,14] = 0.9734851120464462
$data[7,15] = 0.9734851120464463
$data[7,16] = 1050.278080831698
$data[7,17] = 9452.50272748528
$data[7,18] = 0.03876844580709633
$data[7,19] = 0.03876844580709634
$data[8,0] = "sCs"
$data[8,1] = "Col1a2"
$data[8,2] = "Itga11"
$data[8,3] = "sCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 30.939307
$data[8,7] = 92.81792100000001
$data[8,8] = 0.03982438491082609
$data[8,9] = 0.03982438491082609
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.8843623333333334
$data[8,13] = 2.653087
$data[8,14] = 0.02536096827965006
$data[8,15] = 0.02536096827965006
$data[8,16] = 27.36155773023634
$data[8,17] = 246.254019572127
$data[8,18] = 0.001009984962480035
$data[8,19] = 0.001009984962480035

$ws.Range("A2:T10").Value2 = $data
